# Apply crypto price/volume updates from the GitHub Actions data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.376.48'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").Value = '3.792.06'
$ws.Range("E3").Value = '  +1.31%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.07'
$ws.Range("E5").Value = '  +0.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.09'
$ws.Range("E6").Value = '  +0.91%  '

$ws.Range("D7").Value = '3.789.19'
$ws.Range("E7").Value = '  +1.24%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.521'
$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'
$ws.Range("E10").Value = '  +0.79%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.41'
$ws.Range("E11").Value = '  -1.29%  '

$ws.Range("E12").Value = '  -0.36%  '

$ws.Range("E13").Value = '  -1.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.18'
$ws.Range("E14").Value = '  -0.93%  '

$ws.Range("D15").Value = '4.421.22'
$ws.Range("E15").Value = '  +1.11%  '

$ws.Range("D16").Value = '3.792.91'
$ws.Range("E16").Value = '  +1.07%  '

$ws.Range("D17").Value = '68.363.79'
$ws.Range("E17").Value = '  +0.68%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.98'
$ws.Range("E18").Value = '  -1.09%  '

$ws.Range("E19").Value = '  +0.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.97'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.96'
$ws.Range("E21").Value = '  +1.58%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '465.70'
$ws.Range("E22").Value = '  -0.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.699'
$ws.Range("E23").Value = '  -0.60%  '

$ws.Range("E24").Value = '  +9.80%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.00'
$ws.Range("E25").Value = '  +1.08%  '

$ws.Range("E26").Value = '  -2.37%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.93'
$ws.Range("E27").Value = '  -1.73%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.19'
$ws.Range("E28").Value = '  +0.32%  '

$ws.Range("E29").Value = '  +0.10%  '

$ws.Range("E30").Value = '  -0.67%  '

$ws.Range("E31").Value = '  -0.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '30.06'
$ws.Range("E32").Value = '  +0.81%  '

$ws.Range("E33").Value = '  -3.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.15'
$ws.Range("E34").Value = '  +0.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.12%  '

$ws.Range("D36").Value = '3.739.96'
$ws.Range("E36").Value = '  +1.20%  '

$ws.Range("E37").Value = '  -0.92%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.42'
$ws.Range("E38").Value = '  -0.46%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.139'
$ws.Range("E39").Value = '  +0.56%  '

$ws.Range("E40").Value = '  +1.02%  '

$ws.Range("E41").Value = '  +0.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '44.00'
$ws.Range("E44").Value = '  +16.09%  '

$ws.Range("E45").Value = '  -2.13%  '

$ws.Range("E46").Value = '  +3.79%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.42'
$ws.Range("E48").Value = '  -2.27%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '394.63'
$ws.Range("E49").Value = '  +0.35%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '146.66'
$ws.Range("E50").Value = '  +1.97%  '

$ws.Range("D51").Value = '2.798.94'
$ws.Range("E51").Value = '  +4.48%  '
